# trend.xlsx update
# - Result sheet: add helper columns G (=copy of A categories), H (=B/4), I (=C/4)
# - Figure sheet chart: point series at the new H/I columns (quartered values) and
#   G column for categories; add a "Relative Search Interest" title to the value axis
# - Refresh sheet selections / active tab to match the saved state after the edit

$wb = $excel.ActiveWorkbook

$wsResult = $wb.Worksheets.Item("Result")
$wsMulti  = $wb.Worksheets.Item("multiTimeline")
$wsPivot  = $wb.Worksheets.Item("Pivot")
$wsFigure = $wb.Worksheets.Item("Figure")

# ---- Result sheet: build G/H/I helper columns -----------------------------
# Header row: H1/I1 mirror B1/C1 (same shared-string text)
$wsResult.Range("H1").Value2 = $wsResult.Range("B1").Value2
$wsResult.Range("I1").Value2 = $wsResult.Range("C1").Value2

# G2:G21 mirrors A2:A21 (category labels), bold like column A
for ($r = 2; $r -le 21; $r++) {
    $wsResult.Cells.Item($r, 7).Value2 = $wsResult.Cells.Item($r, 1).Value2
    $wsResult.Cells.Item($r, 7).Font.Bold = $true
}

# H2:H21 = B/4, I2:I21 = C/4 (quartered values), using a shared formula like Excel's fill-down
$wsResult.Range("H2").Formula = "=B2/4"
$wsResult.Range("H3:H21").Formula = "=B3/4"
$wsResult.Range("I2").Formula = "=C2/4"
$wsResult.Range("I3:I21").Formula = "=C3/4"

# ---- Figure sheet chart: repoint series to the new helper columns ---------
$chart = $wsFigure.ChartObjects().Item(1).Chart

$ser1 = $chart.SeriesCollection().Item(1)
$ser1.Formula = "=SERIES(Result!`$H`$1,Result!`$G`$2:`$G`$21,Result!`$H`$2:`$H`$21,1)"

$ser2 = $chart.SeriesCollection().Item(2)
$ser2.Formula = "=SERIES(Result!`$I`$1,Result!`$G`$2:`$G`$21,Result!`$I`$2:`$I`$21,2)"

# Add a title to the value (vertical) axis
$valAx = $chart.Axes(2)
$valAx.HasTitle = $true
$valAx.AxisTitle.Text = "Relative Search Interest"

# ---- Selections / active sheet --------------------------------------------
# Touch the non-final sheets first; whichever sheet we Select() on last
# becomes the active tab, so multiTimeline must be last.
$wsPivot.Range("C3:D22").Select() | Out-Null
$wsResult.Range("H2").Select() | Out-Null
$wsFigure.Range("F17").Select() | Out-Null
$wsMulti.Range("A83:C86").Select() | Out-Null
